# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 548 (pushing the existing
# rows 548:615 down to 549:616) and populate the new row with the
# latest week's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 548:615 down one row, creating a blank row 548
$ws.Rows.Item(548).Insert()

# Populate the newly inserted row 548 with this week's record
$ws.Cells.Item(548, 1).Value2 = 8
$ws.Cells.Item(548, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(548, 3).Value2 = "Coquimbo"
$ws.Cells.Item(548, 4).Value2 = 45212
$ws.Cells.Item(548, 5).Value2 = 4
$ws.Cells.Item(548, 6).Value2 = 100114013
$ws.Cells.Item(548, 7).Value2 = "Zanahoria"
$ws.Cells.Item(548, 8).Value2 = "Sin especificar"
$ws.Cells.Item(548, 9).Value2 = "Primera"
$ws.Cells.Item(548, 10).Value2 = 520
$ws.Cells.Item(548, 11).Value2 = 5800
$ws.Cells.Item(548, 12).Value2 = 6000
$ws.Cells.Item(548, 13).Value2 = 5900
$ws.Cells.Item(548, 14).Value2 = "`$/saco 20 kilos"
$ws.Cells.Item(548, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(548, 16).Value2 = 295
$ws.Cells.Item(548, 17).Value2 = 20
$ws.Cells.Item(548, 18).Value2 = "Hortaliza"
